$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to retain exact text representation (values contain
# trailing zeros / thousand-separator dots that Excel would otherwise
# normalise if written as numbers).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '29.301.34'
$ws.Range("E2").Value = '  -0.62%  '
$ws.Range("D3").Value = '1.873.23'
$ws.Range("E3").Value = '  -0.28%  '
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").Value = '0.7088'
$ws.Range("E5").Value = '  -0.82%  '
$ws.Range("D6").Value = '241.83'
$ws.Range("E6").Value = '  -0.17%  '
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D8").Value = '0.3110'
$ws.Range("E8").Value = '  -0.46%  '
$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").Value = '0.07783'
$ws.Range("E9").Value = '  +0.54%  '
$ws.Range("D10").Value = '25.04'
$ws.Range("E10").Value = '  -0.32%  '
$ws.Range("D11").Value = '0.08390'
$ws.Range("E11").Value = '  -0.15%  '
$ws.Range("D12").Value = '1.865.90'
$ws.Range("E12").Value = '  -1.60%  '
$ws.Range("D13").Value = '5.240'
$ws.Range("E13").Value = '  -0.41%  '
$ws.Range("D14").Value = '0.7171'
$ws.Range("E14").Value = '  -0.47%  '
$ws.Range("D15").Value = '91.04'
$ws.Range("E15").Value = '  -0.86%  '
$ws.Range("D16").Value = '29.312.44'
$ws.Range("E16").Value = '  -0.56%  '
$ws.Range("D17").Value = '6.102'
$ws.Range("E17").Value = '  +1.63%  '
$ws.Range("D18").Value = '0.000008308'
$ws.Range("E18").Value = '  +1.18%  '
$ws.Range("D19").Value = '240.05'
$ws.Range("E19").Value = '  -2.00%  '
$ws.Range("B20").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C20").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D20").Value = '2.127.57'
$ws.Range("E20").Value = '  -0.43%  '
$ws.Range("B21").Value = 'Avalanche'
$ws.Range("C21").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D21").Value = '13.21'
$ws.Range("E21").Value = '  -0.31%  '
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  +0.11%  '
$ws.Range("D23").Value = '7.755'
$ws.Range("E23").Value = '  -2.62%  '
$ws.Range("D24").Value = '1.000'
$ws.Range("E24").Value = '  +0.10%  '
$ws.Range("D25").Value = '0.1588'
$ws.Range("E25").Value = '  -2.91%  '
$ws.Range("D26").Value = '162.42'
$ws.Range("E26").Value = '  -1.00%  '
$ws.Range("D27").Value = '9.032'
$ws.Range("E27").Value = '  -0.17%  '
$ws.Range("D28").Value = '18.49'
$ws.Range("E28").Value = '  -1.05%  '
$ws.Range("E29").Value = '  -0.35%  '
$ws.Range("D30").Value = '4.405'
$ws.Range("E30").Value = '  -0.63%  '
$ws.Range("D31").Value = '4.312'
$ws.Range("E31").Value = '  -0.47%  '
$ws.Range("D32").Value = '1.268'
$ws.Range("E32").Value = '  -2.40%  '
$ws.Range("D33").Value = '0.05378'
$ws.Range("E33").Value = '  +2.59%  '
$ws.Range("D34").Value = '1.938'
$ws.Range("E34").Value = '  +0.12%  '
$ws.Range("D35").Value = '0.7503'
$ws.Range("E35").Value = '  -2.96%  '
$ws.Range("D36").Value = '1.176'
$ws.Range("E36").Value = '  -0.29%  '
$ws.Range("D37").Value = '2.683'
$ws.Range("E37").Value = '  +0.23%  '
$ws.Range("D38").Value = '0.01876'
$ws.Range("E38").Value = '  +0.30%  '
$ws.Range("D39").Value = '1.238.48'
$ws.Range("E39").Value = '  +5.22%  '
$ws.Range("D40").Value = '2.731'
$ws.Range("E40").Value = '  +0.27%  '
$ws.Range("D41").Value = '6.520'
$ws.Range("E41").Value = '  +1.22%  '
$ws.Range("D42").Value = '0.8922'
$ws.Range("E42").Value = '  +0.00%  '
$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").Value = '109.10'
$ws.Range("E43").Value = '  +4.51%  '
$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").Value = '72.34'
$ws.Range("E44").Value = '  -2.09%  '
$ws.Range("E45").Value = '  +0.10%  '
$ws.Range("D46").Value = '2.022.79'
$ws.Range("E46").Value = '  -0.22%  '
$ws.Range("E47").Value = '  +7.17%  '
$ws.Range("D48").Value = '0.5199'
$ws.Range("E48").Value = '  +0.04%  '
$ws.Range("E49").Value = '  -1.05%  '
$ws.Range("D50").Value = '9.417'
$ws.Range("E50").Value = '  -0.26%  '
$ws.Range("D51").Value = '0.4339'
$ws.Range("E51").Value = '  +0.26%  '
